$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Snapshot every data row (9..130) across all used columns (A..Z,
#    i.e. 1..26) *before* touching anything, capturing formulas
#    (where present) separately from plain values so we can restore
#    each cell the same way it was stored.
# ------------------------------------------------------------------
$firstRow = 9
$lastRow = 130
$lastCol = 26

$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = New-Object 'object[]' $lastCol
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula()) {
            $rowData[$c - 1] = @{ F = $true; V = $cell.Formula() }
        } else {
            $rowData[$c - 1] = @{ F = $false; V = $cell.Value() }
        }
    }
    $snapshot[$r] = $rowData
}

# ------------------------------------------------------------------
# 2) New row order: for each destination row (9..130) this gives the
#    ORIGINAL row number whose data should end up there.
# ------------------------------------------------------------------
$sourceForRow = @(10,9,11,12,13,15,17,16,18,14,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,36,35,37,38,39,41,48,42,43,44,47,49,40,51,52,50,45,46,97,54,53,60,55,56,57,58,59,62,67,63,64,65,73,74,75,76,68,66,77,79,69,70,81,71,82,72,84,85,86,87,80,78,88,90,91,92,83,94,89,95,96,106,107,98,99,100,101,102,109,110,104,105,61,108,111,112,113,93,114,115,103,116,117,118,119,120,121,122,123,124,125,126,127,128,129,130)

# ------------------------------------------------------------------
# 3) Write the snapshot back out in the new order.
# ------------------------------------------------------------------
for ($i = 0; $i -lt $sourceForRow.Length; $i++) {
    $destRow = $firstRow + $i
    $srcRow = $sourceForRow[$i]
    $rowData = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $entry = $rowData[$c - 1]
        $cell = $ws.Cells.Item($destRow, $c)
        if ($entry.F) {
            $cell.Formula = $entry.V
        } else {
            $cell.Value = $entry.V
        }
    }
}

# ------------------------------------------------------------------
# 4) "Förändrad" (column C) is bumped to 46072 for every data row.
# ------------------------------------------------------------------
for ($r = 2; $r -le 130; $r++) {
    $ws.Cells.Item($r, 3).Value = 46072
}
